# Update countries & provincias Spain
#
# This script reproduces, cell-by-cell, the data refresh captured in the
# target diff:
#  - the "Datos actualizados" timestamp in A1 moves from 12:46 to 13:16
#  - "A Coruña" (row 12) gets an updated case/death count
#  - "Caceres" jumps up the ranking to row 20 (pushing Salamanca..Zaragoza
#    down by one row each, rows 21-28) with a brand new set of numbers
#  - "Badajoz" jumps up the ranking to row 34 (pushing Jaen..Aragon down by
#    one row each, rows 35-39) with a brand new set of numbers

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 24 de Marzo de 2020 a las 13:16"

function Set-Row($row, $name, $total, $active, $recovered, $deaths) {
    $ws.Range("A$row").Value = $name
    $ws.Range("B$row").Value = $total
    $ws.Range("C$row").Value = $active
    $ws.Range("D$row").Value = $recovered
    $ws.Range("E$row").Value = $deaths
}

Set-Row 12 "A Coruña" 635 19 621 14

Set-Row 20 "Caceres" 419 3 392 24
Set-Row 21 "Salamanca" 404 22 358 24
Set-Row 22 "Murcia" 385 1 381 3
Set-Row 23 "Gipuzkoa/Guipuzcoa" 380 283 365 15
Set-Row 24 "Granada" 374 0 357 17
Set-Row 25 "Sevilla" 351 1 345 5
Set-Row 26 "Valladolid" 349 17 318 14
Set-Row 27 "Burgos" 336 29 289 18
Set-Row 28 "Zaragoza" 329 0 315 14

Set-Row 34 "Badajoz" 217 5 209 3
Set-Row 35 "Jaen" 215 0 210 5
Set-Row 36 "Mallorca" 210 18 194 12
Set-Row 37 "Cordoba" 191 0 187 4
Set-Row 38 "Cadiz" 178 0 175 3
Set-Row 39 "Aragon" 174 0 163 11
